$d = $word.ActiveDocument

# 1) Patient name: OTAVIO RAMOS DE ALMEIDA -> ARNALDO BALDOW
$d.Content.Find.Execute("OTAVIO RAMOS DE ALMEIDA", $true, $false, $false, $false, $false, $true, 1, $false, "ARNALDO BALDOW", 2) | Out-Null

# 2) Birth date: 15/02/1988 -> 05/09/1952
$d.Content.Find.Execute("15/02/1988", $true, $false, $false, $false, $false, $true, 1, $false, "05/09/1952", 2) | Out-Null

# 3) Record number (Porntuario): 26294 -> 126577
$d.Content.Find.Execute("26294", $true, $false, $false, $false, $false, $true, 1, $false, "126577", 2) | Out-Null

# 4) Mother's name: MARIA IRACY RAMOS DOS SANTOS -> ALMIRA MOREIRA BALDOW
$d.Content.Find.Execute("MARIA IRACY RAMOS DOS SANTOS", $true, $false, $false, $false, $false, $true, 1, $false, "ALMIRA MOREIRA BALDOW", 2) | Out-Null

# 5) Record date: 03/01/2019 -> 11/05/2018
$d.Content.Find.Execute("03/01/2019", $true, $false, $false, $false, $false, $true, 1, $false, "11/05/2018", 2) | Out-Null

# 6) Clear the urology note text at the end of the "Evolucao Consultorio"
#    paragraph, leaving the trailing run present but with no text (the
#    printed note text is removed, per "ESCONDER BOTAO IMPRIMIR").
$noteText = "# UROLOGIA`nSOLICITO RETIRADA DE DUPLO J`n30 PO DE URETERO + DUPLO J"
$full = $d.Content.Text
$startIdx = $full.IndexOf($noteText)
if ($startIdx -ge 0) {
    $noteRange = $d.Range($startIdx, $startIdx + $noteText.Length)
    $noteRange.Text = ""

    # Re-insert an empty run at that same spot so a (now empty) run remains
    # in place of the note, matching the trailing empty <w:r> left behind.
    $placeholder = $d.Range($startIdx, $startIdx)
    $placeholder.InsertAfter("")
}
